# Project "Sample Project" save: update the "Rules" sheet, cell B11.
# B11 currently holds the shared string "R40"; it must become the
# (text) string "1" while keeping its existing cell style (s="23",
# t="s") untouched.
#
# A plain `$ws.Range("B11").Value = "1"` (or `.Value2 = "1"`) gets
# auto-coerced to the NUMBER 1 by Excel's type inference, which both
# changes the cell's type (t="s" -> no "t" attr) and drops it out of
# the shared-string table. Forcing the format to "@"/Text first also
# works, but it mutates the cell's style (a new xf/numFmt gets
# allocated), which is not part of this edit.
#
# Instead, write it as a text formula (="1") and then convert that
# formula to its static value in place via Copy / PasteSpecial
# (values only) -- this is exactly what "Copy, Paste Values" in the
# Excel UI does, and it preserves the original number format / style
# index while leaving the cell holding a literal text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
